$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 374, shifting existing rows 374-456 down to 375-457.
$ws.Rows.Item(374).Insert()

# Populate the new row 374 with its data.
$ws.Range("A374").Value = 7
$ws.Range("B374").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C374").Value = "Ñuble"
$ws.Range("D374").Value = 45275
$ws.Range("E374").Value = 16
$ws.Range("F374").Value = "Fruta"
$ws.Range("G374").Value = 100103
$ws.Range("H374").Value = "Frutos de hueso (carozo)"
$ws.Range("I374").Value = 100103004
$ws.Range("J374").Value = "Durazno"
$ws.Range("K374").Value = "Rich Lady"
$ws.Range("L374").Value = "Primera"
$ws.Range("M374").Value = 60
$ws.Range("N374").Value = 15000
$ws.Range("O374").Value = 15000
$ws.Range("P374").Value = 15000
$ws.Range("Q374").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R374").Value = "Región de O'Higgins"
$ws.Range("S374").Value = 938
$ws.Range("T374").Value = 16
